$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 55.21
$ws.Range("I15").Value = 55.21
$ws.Range("K15").Value = 165.63
$ws.Range("M15").Value = 3.370000000000005

$ws.Range("H55").Value = 234.5
$ws.Range("I55").Value = 138
$ws.Range("J55").Value = 303.42856
$ws.Range("K55").Value = 138
$ws.Range("L55").Value = 303.42856
$ws.Range("M55").Value = 76
$ws.Range("N55").Value = -731.4285600000001

$ws.Range("H62").Value = 15321.723
$ws.Range("I62").Value = 22799.1
$ws.Range("J62").Value = 5975
$ws.Range("K62").Value = 22799.1
$ws.Range("L62").Value = 5975
$ws.Range("M62").Value = -22175.1
$ws.Range("N62").Value = -7223

$ws.Range("H65").Value = 15321.723
$ws.Range("I65").Value = 22799.1
$ws.Range("J65").Value = 5975
$ws.Range("K65").Value = 113995.5
$ws.Range("L65").Value = 29875
$ws.Range("M65").Value = -110875.5
$ws.Range("N65").Value = -36115

$ws.Range("H98").Value = 661.1053000000001
$ws.Range("I98").Value = 675.6111
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 675.6111
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = 822.3889
$ws.Range("N98").Value = -3396

$ws.Range("H116").Value = 36084.562
$ws.Range("I116").Value = 79036.86
$ws.Range("J116").Value = 2677.2222
$ws.Range("K116").Value = 79036.86
$ws.Range("L116").Value = 2677.2222
$ws.Range("M116").Value = -75594.86
$ws.Range("N116").Value = -9561.2222

$ws.Range("H122").Value = 661.1053000000001
$ws.Range("I122").Value = 675.6111
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 2026.8333
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = 423.1667000000002
$ws.Range("N122").Value = -6100

$ws.Range("H132").Value = 2270.3672
$ws.Range("I132").Value = 1204.1111
$ws.Range("J132").Value = 6468.75
$ws.Range("K132").Value = 3612.3333
$ws.Range("L132").Value = 19406.25
$ws.Range("M132").Value = -1082.3333
$ws.Range("N132").Value = -24466.25

$ws.Range("H138").Value = 1528.9572
$ws.Range("I138").Value = 1020.91895
$ws.Range("J138").Value = 2098.5757
$ws.Range("K138").Value = 3062.75685
$ws.Range("L138").Value = 6295.7271
$ws.Range("M138").Value = 2077.24315
$ws.Range("N138").Value = -16575.7271

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3795.8462
$ws.Range("I2").Value = 2567.818
$ws.Range("J2").Value = 10550
$ws.Range("K2").Value = 2567.818
$ws.Range("L2").Value = 10550
$ws.Range("M2").Value = -2454.818
$ws.Range("N2").Value = -10776

$ws.Range("H32").Value = 1718.84
$ws.Range("I32").Value = 1468
$ws.Range("J32").Value = 4603.5
$ws.Range("K32").Value = 1468
$ws.Range("L32").Value = 4603.5
$ws.Range("M32").Value = -1181
$ws.Range("N32").Value = -5177.5

$ws.Range("H61").Value = 302142
$ws.Range("I61").Value = 221350.27
$ws.Range("J61").Value = 504121.4
$ws.Range("K61").Value = 221350.27
$ws.Range("L61").Value = 504121.4
$ws.Range("M61").Value = -221138.27
$ws.Range("N61").Value = -504545.4

$ws.Range("H116").Value = 3795.8462
$ws.Range("I116").Value = 2567.818
$ws.Range("J116").Value = 10550
$ws.Range("K116").Value = 2567.818
$ws.Range("L116").Value = 10550
$ws.Range("M116").Value = -273.8180000000002
$ws.Range("N116").Value = -15138

$ws.Range("H122").Value = 3862.0222
$ws.Range("I122").Value = 3867.0278
$ws.Range("J122").Value = 3842
$ws.Range("K122").Value = 11601.0834
$ws.Range("L122").Value = 11526
$ws.Range("M122").Value = -9151.0834
$ws.Range("N122").Value = -16426

$ws.Range("H132").Value = 2149.1704
$ws.Range("I132").Value = 1869.2084
$ws.Range("J132").Value = 3409
$ws.Range("K132").Value = 5607.6252
$ws.Range("L132").Value = 10227
$ws.Range("M132").Value = -3077.6252
$ws.Range("N132").Value = -15287

$ws.Range("H136").Value = 302142
$ws.Range("I136").Value = 221350.27
$ws.Range("J136").Value = 504121.4
$ws.Range("K136").Value = 664050.8099999999
$ws.Range("L136").Value = 1512364.2
$ws.Range("M136").Value = -661500.8099999999
$ws.Range("N136").Value = -1517464.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3795.8462
$ws.Range("I3").Value = 2567.818
$ws.Range("J3").Value = 10550
$ws.Range("K3").Value = 2567.818
$ws.Range("L3").Value = 10550
$ws.Range("M3").Value = -2453.818
$ws.Range("N3").Value = -10778

$ws.Range("H80").Value = 127.388885
$ws.Range("I80").Value = 97.71429000000001
$ws.Range("J80").Value = 146.27272
$ws.Range("K80").Value = 97.71429000000001
$ws.Range("L80").Value = 146.27272
$ws.Range("M80").Value = 900.28571
$ws.Range("N80").Value = -2142.27272

$ws.Range("H83").Value = 127.388885
$ws.Range("I83").Value = 97.71429000000001
$ws.Range("J83").Value = 146.27272
$ws.Range("K83").Value = 488.57145
$ws.Range("L83").Value = 731.3635999999999
$ws.Range("M83").Value = 4503.42855
$ws.Range("N83").Value = -10715.3636

$ws.Range("H99").Value = 4949813.5
$ws.Range("I99").Value = 1675748.4
$ws.Range("J99").Value = 17500396
$ws.Range("K99").Value = 1675748.4
$ws.Range("L99").Value = 17500396
$ws.Range("M99").Value = -1674250.4
$ws.Range("N99").Value = -17503392

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 685.5833
$ws.Range("I16").Value = 628.8
$ws.Range("J16").Value = 726.1429000000001
$ws.Range("K16").Value = 628.8
$ws.Range("L16").Value = 726.1429000000001
$ws.Range("M16").Value = -341.8
$ws.Range("N16").Value = -1300.1429

$ws.Range("H31").Value = 2092.4531
$ws.Range("I31").Value = 1221.902
$ws.Range("K31").Value = 1221.902
$ws.Range("M31").Value = -926.902

$ws.Range("H34").Value = 2092.4531
$ws.Range("I34").Value = 1221.902
$ws.Range("K34").Value = 1221.902
$ws.Range("M34").Value = -1019.902

$ws.Range("H58").Value = 3628.7727
$ws.Range("I58").Value = 4433.5
$ws.Range("J58").Value = 2220.5
$ws.Range("K58").Value = 4433.5
$ws.Range("L58").Value = 2220.5
$ws.Range("M58").Value = -4230.5
$ws.Range("N58").Value = -2626.5

$ws.Range("H94").Value = 6816.2144
$ws.Range("I94").Value = 925
$ws.Range("K94").Value = 925
$ws.Range("M94").Value = -474

$ws.Range("H105").Value = 834.8823
$ws.Range("I105").Value = 500.9091
$ws.Range("J105").Value = 1447.1666
$ws.Range("K105").Value = 500.9091
$ws.Range("L105").Value = 1447.1666
$ws.Range("M105").Value = 1246.0909
$ws.Range("N105").Value = -4941.1666

$ws.Range("H113").Value = 685.5833
$ws.Range("I113").Value = 628.8
$ws.Range("J113").Value = 726.1429000000001
$ws.Range("K113").Value = 628.8
$ws.Range("L113").Value = 726.1429000000001
$ws.Range("M113").Value = 1541.2
$ws.Range("N113").Value = -5066.1429

$ws.Range("H132").Value = 1735.92
$ws.Range("I132").Value = 1105.1082
$ws.Range("J132").Value = 3531.3076
$ws.Range("K132").Value = 3315.3246
$ws.Range("L132").Value = 10593.9228
$ws.Range("M132").Value = -785.3245999999999
$ws.Range("N132").Value = -15653.9228

$ws.Range("H136").Value = 3628.7727
$ws.Range("I136").Value = 4433.5
$ws.Range("J136").Value = 2220.5
$ws.Range("K136").Value = 13300.5
$ws.Range("L136").Value = 6661.5
$ws.Range("M136").Value = -10750.5
$ws.Range("N136").Value = -11761.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 411.72223
$ws.Range("I5").Value = 363.1875
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 1089.5625
$ws.Range("L5").Value = 2400
$ws.Range("M5").Value = -977.5625
$ws.Range("N5").Value = -2624

$ws.Range("H122").Value = 750.3
$ws.Range("I122").Value = 489.0909
$ws.Range("J122").Value = 901.5263
$ws.Range("K122").Value = 4401.8181
$ws.Range("L122").Value = 8113.736699999999
$ws.Range("M122").Value = -1951.8181
$ws.Range("N122").Value = -13013.7367

$ws.Range("H135").Value = 411.72223
$ws.Range("I135").Value = 363.1875
$ws.Range("J135").Value = 800
$ws.Range("K135").Value = 3268.6875
$ws.Range("L135").Value = 7200
$ws.Range("M135").Value = -733.6875
$ws.Range("N135").Value = -12270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 985.125
$ws.Range("I122").Value = 984.2
$ws.Range("J122").Value = 999
$ws.Range("K122").Value = 2952.6
$ws.Range("L122").Value = 2997
$ws.Range("M122").Value = -502.6000000000004
$ws.Range("N122").Value = -7897

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1537.2715
$ws.Range("I132").Value = 1267.9138
$ws.Range("K132").Value = 3803.7414
$ws.Range("M132").Value = -1273.7414
